$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sv1")

# Delete the entire row that contains station "OKKR" (data point #199).
# That row is the physical worksheet row 200 (header is row 1, first data
# row is row 2, so data point 199 sits on row 201... but OKKR is actually
# the row whose sequence number in column A is 199, i.e. worksheet row 200).
$ws.Rows.Item(200).Delete()
